# Apply the commit: append newly scraped job listing as the new top row
# on sheet 1 (ランサーズ) and append a new statistics row on sheet 2 (統計).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: ランサーズ (job listings) ---------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a brand-new row above the current row 2, shifting every existing
# listing (and its hyperlink) down by one row.
$ws1.Rows.Item(2).Insert()

# Fill in the newly scraped listing.
$ws1.Range("A2").Value = "2025-09-03 01:40:17"
$ws1.Range("B2").Value = "Google口コミ促進ツールの制作"
$ws1.Range("C2").Value = "システム開発"
$ws1.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws1.Range("E2").Value = "期限情報なし"
$ws1.Range("F2").Value = "https://www.lancers.jp/work/detail/5385483"
$ws1.Range("G2").Value = 80
$ws1.Range("H2").Value = "◆ツール"

# Recreate the hyperlink on the URL cell, then restore the normal
# "Hyperlink" cell style used by every other row in column F.
$ws1.Hyperlinks.Add($ws1.Range("F2"), "https://www.lancers.jp/work/detail/5385483")
$ws1.Range("F2").Style = $ws1.Range("F3").Style

# --- Sheet 2: 統計 (stats log) ------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# This log is append-only; add the new summary row at the bottom.
$ws2.Range("A56").Value = "2025-09-03T01:40:17.435645"
$ws2.Range("B56").Value = 23
$ws2.Range("C56").Value = "全案件リスト"
$ws2.Range("D56").Value = 73.90000000000001
$ws2.Range("E56").Value = 6
$ws2.Range("F56").Value = 10
$ws2.Range("G56").Value = 23
